$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 111111840
$ws.Range("I42").Value = 200000510
$ws.Range("K42").Value = 600001530
$ws.Range("M42").Value = -600001300
$ws.Range("H53").Value = 372.1613
$ws.Range("J53").Value = 581.1
$ws.Range("L53").Value = 581.1
$ws.Range("N53").Value = -1855.1
$ws.Range("H74").Value = 4713.4287
$ws.Range("I74").Value = 4713.4287
$ws.Range("K74").Value = 4713.4287
$ws.Range("M74").Value = -3777.4287
$ws.Range("H77").Value = 4713.4287
$ws.Range("I77").Value = 4713.4287
$ws.Range("K77").Value = 23567.1435
$ws.Range("M77").Value = -18887.1435
$ws.Range("H132").Value = 1637.0571
$ws.Range("I132").Value = 1269.129
$ws.Range("K132").Value = 3807.387
$ws.Range("M132").Value = -1277.387
$ws.Range("H137").Value = 22499.8
$ws.Range("I137").Value = 100002
$ws.Range("J137").Value = 3124.25
$ws.Range("K137").Value = 300006
$ws.Range("L137").Value = 9372.75
$ws.Range("M137").Value = -297456
$ws.Range("N137").Value = -14472.75
$ws.Range("H138").Value = 2451.5532
$ws.Range("I138").Value = 1231.4445
$ws.Range("J138").Value = 2740.5264
$ws.Range("K138").Value = 3694.3335
$ws.Range("L138").Value = 8221.5792
$ws.Range("M138").Value = 1445.6665
$ws.Range("N138").Value = -18501.5792
$ws.Range("H141").Value = 8191.65
$ws.Range("I141").Value = 7308.6665
$ws.Range("K141").Value = 21925.9995
$ws.Range("M141").Value = -16745.9995

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9757.75
$ws.Range("I32").Value = 10099.363
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 10099.363
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -9812.362999999999
$ws.Range("N32").Value = -6574
$ws.Range("H38").Value = 5369.5713
$ws.Range("I38").Value = 5277.6665
$ws.Range("J38").Value = 5921
$ws.Range("K38").Value = 5277.6665
$ws.Range("L38").Value = 5921
$ws.Range("M38").Value = -4810.6665
$ws.Range("N38").Value = -6855
$ws.Range("H45").Value = 1764.0476
$ws.Range("I45").Value = 1277.0834
$ws.Range("J45").Value = 2413.3333
$ws.Range("K45").Value = 1277.0834
$ws.Range("L45").Value = 2413.3333
$ws.Range("M45").Value = -900.0834
$ws.Range("N45").Value = -3167.3333
$ws.Range("H61").Value = 2705.8235
$ws.Range("I61").Value = 2656.1875
$ws.Range("K61").Value = 2656.1875
$ws.Range("M61").Value = -2444.1875
$ws.Range("H63").Value = 4503.6
$ws.Range("I63").Value = 2656.4443
$ws.Range("J63").Value = 7274.3335
$ws.Range("K63").Value = 2656.4443
$ws.Range("L63").Value = 7274.3335
$ws.Range("M63").Value = -1970.4443
$ws.Range("N63").Value = -8646.333500000001
$ws.Range("H66").Value = 4503.6
$ws.Range("I66").Value = 2656.4443
$ws.Range("J66").Value = 7274.3335
$ws.Range("K66").Value = 13282.2215
$ws.Range("L66").Value = 36371.6675
$ws.Range("M66").Value = -9850.2215
$ws.Range("N66").Value = -43235.6675
$ws.Range("H74").Value = 1959.0588
$ws.Range("I74").Value = 1987.75
$ws.Range("K74").Value = 1987.75
$ws.Range("M74").Value = -1113.75
$ws.Range("H77").Value = 1959.0588
$ws.Range("I77").Value = 1987.75
$ws.Range("K77").Value = 9938.75
$ws.Range("M77").Value = -5570.75
$ws.Range("H122").Value = 1499.683
$ws.Range("I122").Value = 1218.3939
$ws.Range("J122").Value = 2660
$ws.Range("K122").Value = 3655.1817
$ws.Range("L122").Value = 7980
$ws.Range("M122").Value = -1205.1817
$ws.Range("N122").Value = -12880
$ws.Range("H124").Value = 28994.5
$ws.Range("J124").Value = 28994.5
$ws.Range("L124").Value = 28994.5
$ws.Range("N124").Value = -38814.5
$ws.Range("H125").Value = 64999.5
$ws.Range("J125").Value = 64999.5
$ws.Range("L125").Value = 64999.5
$ws.Range("N125").Value = -74839.5
$ws.Range("H136").Value = 2705.8235
$ws.Range("I136").Value = 2656.1875
$ws.Range("K136").Value = 7968.5625
$ws.Range("M136").Value = -5418.5625

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 51303.125
$ws.Range("J132").Value = 51303.125
$ws.Range("L132").Value = 51303.125
$ws.Range("N132").Value = -61423.125
$ws.Range("H138").Value = 50780
$ws.Range("J138").Value = 50780
$ws.Range("L138").Value = 50780
$ws.Range("N138").Value = -61060

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5460.1333
$ws.Range("I31").Value = 1809.5454
$ws.Range("K31").Value = 1809.5454
$ws.Range("M31").Value = -1514.5454
$ws.Range("H34").Value = 5460.1333
$ws.Range("I34").Value = 1809.5454
$ws.Range("K34").Value = 1809.5454
$ws.Range("M34").Value = -1607.5454
$ws.Range("H43").Value = 34919.75
$ws.Range("J43").Value = 34919.75
$ws.Range("L43").Value = 34919.75
$ws.Range("N43").Value = -35287.75
$ws.Range("H101").Value = 34919.75
$ws.Range("J101").Value = 34919.75
$ws.Range("L101").Value = 34919.75
$ws.Range("N101").Value = -41409.75
$ws.Range("H122").Value = 1899.2858
$ws.Range("I122").Value = 1859
$ws.Range("K122").Value = 5577
$ws.Range("M122").Value = -3127
$ws.Range("H132").Value = 2080.3794
$ws.Range("I132").Value = 2047.0358
$ws.Range("K132").Value = 6141.107400000001
$ws.Range("M132").Value = -3611.107400000001

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 985.64703
$ws.Range("I5").Value = 450.6
$ws.Range("J5").Value = 1750
$ws.Range("K5").Value = 1351.8
$ws.Range("L5").Value = 5250
$ws.Range("M5").Value = -1239.8
$ws.Range("N5").Value = -5474
$ws.Range("H33").Value = 142.16667
$ws.Range("J33").Value = 200.66667
$ws.Range("L33").Value = 1204.00002
$ws.Range("N33").Value = -1770.00002
$ws.Range("H109").Value = 53757.43
$ws.Range("I109").Value = 61761.766
$ws.Range("K109").Value = 185285.298
$ws.Range("M109").Value = -184245.298
$ws.Range("H113").Value = 1300.6666
$ws.Range("J113").Value = 1417.4667
$ws.Range("L113").Value = 4252.4001
$ws.Range("N113").Value = -8592.400099999999
$ws.Range("H133").Value = 1472.2858
$ws.Range("I133").Value = 1472.2858
$ws.Range("K133").Value = 4416.857400000001
$ws.Range("M133").Value = 643.1425999999992
$ws.Range("H135").Value = 985.64703
$ws.Range("I135").Value = 450.6
$ws.Range("J135").Value = 1750
$ws.Range("K135").Value = 4055.4
$ws.Range("L135").Value = 15750
$ws.Range("M135").Value = -1520.4
$ws.Range("N135").Value = -20820

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21998.334
$ws.Range("I70").Value = 17999
$ws.Range("K70").Value = 17999
$ws.Range("M70").Value = -17729
$ws.Range("H73").Value = 21998.334
$ws.Range("I73").Value = 17999
$ws.Range("K73").Value = 17999
$ws.Range("M73").Value = -17063
$ws.Range("H80").Value = 8078.7144
$ws.Range("I80").Value = 12186.728
$ws.Range("J80").Value = 3559.9
$ws.Range("K80").Value = 12186.728
$ws.Range("L80").Value = 3559.9
$ws.Range("M80").Value = -11188.728
$ws.Range("N80").Value = -5555.9
$ws.Range("H83").Value = 8078.7144
$ws.Range("I83").Value = 12186.728
$ws.Range("J83").Value = 3559.9
$ws.Range("K83").Value = 60933.64
$ws.Range("L83").Value = 17799.5
$ws.Range("M83").Value = -55941.64
$ws.Range("N83").Value = -27783.5

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2111.111
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1888
$ws.Range("H7").Value = 4515.4443
$ws.Range("J7").Value = 3875
$ws.Range("L7").Value = 3875
$ws.Range("N7").Value = -4099
$ws.Range("H15").Value = 2111.111
$ws.Range("I15").Value = 2000
$ws.Range("K15").Value = 2000
$ws.Range("M15").Value = -1830
$ws.Range("H16").Value = 664.7273
$ws.Range("I16").Value = 664.7273
$ws.Range("K16").Value = 664.7273
$ws.Range("M16").Value = -494.7273
$ws.Range("H24").Value = 4500
$ws.Range("I24").Value = 4500
$ws.Range("K24").Value = 4500
$ws.Range("M24").Value = -4157
$ws.Range("H59").Value = 21000
$ws.Range("J59").Value = 21000
$ws.Range("L59").Value = 21000
$ws.Range("N59").Value = -22308
$ws.Range("H126").Value = 4515.4443
$ws.Range("J126").Value = 3875
$ws.Range("L126").Value = 11625
$ws.Range("N126").Value = -16565
$ws.Range("H132").Value = 5010.6313
$ws.Range("J132").Value = 4945.727
$ws.Range("L132").Value = 14837.181
$ws.Range("N132").Value = -19897.181

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2506000
$ws.Range("I3").Value = 5000000
$ws.Range("J3").Value = 12000
$ws.Range("K3").Value = 5000000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -4999886
$ws.Range("N3").Value = -12228
$ws.Range("H70").Value = 16666.666
$ws.Range("I70").Value = 16666.666
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 16666.666
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -16351.666
$ws.Range("H73").Value = 16666.666
$ws.Range("I73").Value = 16666.666
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 16666.666
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -15574.666
$ws.Range("H122").Value = 7058.1333
$ws.Range("J122").Value = 9797.9
$ws.Range("L122").Value = 29393.7
$ws.Range("N122").Value = -34293.7
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()
